# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# worksheet to the latest scraped snapshot. Each entry below carries the
# new Price (D) / Volume(1h) (E) text for its row; a $null column means
# that cell did not change in this refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "28.278.83"; E = "  +0.22%  " },
    @{ Row = 3; D = "1.858.45"; E = "  -0.61%  " },
    @{ Row = 4; D = "1.000"; E = "  +0.05%  " },
    @{ Row = 5; D = "329.86"; E = "  -2.21%  " },
    @{ Row = 6; D = $null; E = "  +0.04%  " },
    @{ Row = 7; D = "0.4555"; E = "  -3.16%  " },
    @{ Row = 8; D = "0.3935"; E = "  +0.18%  " },
    @{ Row = 9; D = "47.56"; E = "  +1.05%  " },
    @{ Row = 10; D = "0.07817"; E = "  -2.12%  " },
    @{ Row = 11; D = "0.9821"; E = "  -2.64%  " },
    @{ Row = 12; D = $null; E = "  -2.09%  " },
    @{ Row = 13; D = "1.852.93"; E = "  +1.19%  " },
    @{ Row = 14; D = "5.805"; E = "  -3.18%  " },
    @{ Row = 15; D = "6.956"; E = "  -4.41%  " },
    @{ Row = 16; D = "1.003"; E = "  +0.20%  " },
    @{ Row = 17; D = "88.07"; E = "  -3.59%  " },
    @{ Row = 18; D = "0.06536"; E = "  -0.96%  " },
    @{ Row = 19; D = "0.00001018"; E = "  -1.91%  " },
    @{ Row = 20; D = "17.07"; E = "  -3.77%  " },
    @{ Row = 21; D = "1.002"; E = "  +0.13%  " },
    @{ Row = 22; D = "28.238.64"; E = "  +0.10%  " },
    @{ Row = 23; D = "5.286"; E = "  -2.91%  " },
    @{ Row = 24; D = "10.81"; E = "  -2.05%  " },
    @{ Row = 25; D = "2.255"; E = "  -1.55%  " },
    @{ Row = 26; D = "2.076.88"; E = "  +1.13%  " },
    @{ Row = 27; D = "156.80"; E = "  -1.44%  " },
    @{ Row = 28; D = $null; E = "  -3.77%  " },
    @{ Row = 29; D = "2.048"; E = "  -4.10%  " },
    @{ Row = 30; D = "5.274"; E = "  -3.83%  " },
    @{ Row = 31; D = "116.40"; E = "  -2.83%  " },
    @{ Row = 32; D = $null; E = "  -2.21%  " },
    @{ Row = 33; D = "0.9365"; E = "  -4.09%  " },
    @{ Row = 34; D = "3.604"; E = "  +0.82%  " },
    @{ Row = 35; D = "1.378"; E = "  -0.02%  " },
    @{ Row = 36; D = "5.181"; E = "  -3.13%  " },
    @{ Row = 37; D = "0.06012"; E = "  -1.54%  " },
    @{ Row = 38; D = "0.02206"; E = "  -2.77%  " },
    @{ Row = 39; D = "8.255"; E = "  -2.23%  " },
    @{ Row = 40; D = "1.160"; E = "  -1.37%  " },
    @{ Row = 41; D = "1.001"; E = "  +0.14%  " },
    @{ Row = 42; D = "0.5673"; E = "  -5.21%  " },
    @{ Row = 43; D = "9.971"; E = "  -3.61%  " },
    @{ Row = 44; D = "0.1790"; E = "  -5.05%  " },
    @{ Row = 45; D = "1.245"; E = "  -2.43%  " },
    @{ Row = 46; D = "2.296"; E = "  +16.56%  " },
    @{ Row = 47; D = "0.5390"; E = "  -4.33%  " },
    @{ Row = 48; D = "11.79"; E = "  -3.58%  " },
    @{ Row = 49; D = "0.07164"; E = "  +4.33%  " },
    @{ Row = 50; D = "1.862"; E = "  -5.90%  " },
    @{ Row = 51; D = "109.84"; E = "  -1.60%  " }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($row, 4)   # column D = Price
        # Many prices ("1.000", "329.86", "0.07817", ...) look numeric;
        # force text formatting first so Excel stores the literal digits
        # instead of silently parsing/rounding them as a number.
        if ($u.D -match '^-?[0-9]*\.?[0-9]+$') {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $u.D
    }

    if ($null -ne $u.E) {
        $ws.Cells.Item($row, 5).Value = $u.E   # column E = Volume(1h)
    }
}
